$d = $word.ActiveDocument

# The document currently ends with an empty paragraph that holds the
# "_GoBack" bookmark (last paragraph in the body, right before sectPr).
# We replace that whole paragraph with: a blank paragraph, the bookmark
# paragraph now carrying "zad.3" text, and three more paragraphs of
# instructions -- matching the target OOXML exactly via InsertXML.
$target = $d.Paragraphs.Last.Range

$xml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="0092067C" w:rsidRDefault="0092067C">
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
      </w:pPr>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve">zad.3 Napisać funkcje </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>ThrowSomething</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve">, </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>ktora</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> na podstawie przekazanego do niej typu </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>zwroci</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>wyjatek</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> tego typu wraz z </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>jakas</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>wartoscia</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve">Ma </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>wykrywac</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> pierwsze 4 typy z </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>enum</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>Type</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve">, natomiast w innym wypadku ma </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>wyrzucac</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> tekst "Nieznany typ!".</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
      </w:pPr>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>Napisac</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> odpowiednio </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>zagniezdzone</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> instrukcje </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>try</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>/</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>catch</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> w funkcji </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>ExceptionHandler</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>,</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:autoSpaceDE w:val="0"/>
        <w:autoSpaceDN w:val="0"/>
        <w:adjustRightInd w:val="0"/>
        <w:spacing w:after="0" w:line="240" w:lineRule="auto"/>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve">tak by najpierw </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>lapala</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> znane typy, a dopiero w </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>ostatecznosci</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> nieznane typy.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="000000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve">W </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>kazdym</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> bloku </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>catch</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> ma </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>zostac</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> wypisana nazwa typu oraz </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>wartosc</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t xml:space="preserve"> rzuconego </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>wyjatku</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Consolas" w:hAnsi="Consolas" w:cs="Consolas"/>
          <w:color w:val="008000"/>
          <w:sz w:val="19"/>
          <w:szCs w:val="19"/>
        </w:rPr>
        <w:t>.</w:t>
      </w:r>
    </w:p>
"@

$target.InsertXML($xml)
